$d = $word.ActiveDocument
$p = $d.Paragraphs.Item(1)

# Add paragraph border spacing (top/left/bottom/right "space" = 5, no line)
$p.Format.Borders.DistanceFromTop = 5
$p.Format.Borders.DistanceFromLeft = 5
$p.Format.Borders.DistanceFromBottom = 5
$p.Format.Borders.DistanceFromRight = 5

# Update left indent from 120 twips (6pt) to 225 twips (11.25pt)
$p.Format.LeftIndent = 11.25

# Replace the paragraph's text (everything up to, but not including, the
# paragraph mark) with the new ID, collapsing both runs (the ID run and the
# trailing-space run) into a single run that keeps the first run's
# formatting.
$r = $d.Range($p.Range.Start, $p.Range.End - 1)
$r.Text = "**ID__AFFARS_5337_7401__ID**"
